$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM values for rows 2-7 (Neutrophils / Il1f10 / Il1r1 pairs)
$data = @{
    2 = @{ G = 2.751441333333334; H = 8.254324; M = 16.087096; N = 48.261288; O = 0.1263055268415452; P = 0.1263055268415452; Q = 44.26270086770133; R = 398.364307809312; S = 0.1263055268415452; T = 0.1263055268415452 }
    3 = @{ G = 2.751441333333334; H = 8.254324; O = 0.7490048915888087; P = 0.7490048915888088; Q = 262.4824130335334; R = 2362.3417173018; S = 0.7490048915888087; T = 0.7490048915888088 }
    4 = @{ G = 2.751441333333334; H = 8.254324; M = 0.5200936666666667; N = 1.560281; O = 0.004083440825819921; P = 0.004083440825819921; Q = 1.431007211671556; R = 12.879064905044; S = 0.004083440825819921; T = 0.004083440825819921 }
    5 = @{ G = 2.751441333333334; H = 8.254324; M = 14.15205133333333; N = 42.456154; O = 0.1111128011883101; P = 0.1111128011883101; Q = 38.93853898998844; R = 350.446850909896; S = 0.1111128011883101; T = 0.1111128011883101 }
    6 = @{ G = 2.751441333333334; H = 8.254324; M = 0.794831; N = 2.384493; O = 0.006240501592393819; P = 0.006240501592393819; Q = 2.186930866414667; R = 19.682377797732; S = 0.006240501592393819; T = 0.006240501592393819 }
    7 = @{ G = 2.751441333333334; H = 8.254324; M = 0.4143026666666667; N = 1.242908; O = 0.003252837963122146; P = 0.003252837963122146; Q = 1.139929481576889; R = 10.259365334192; S = 0.003252837963122146; T = 0.003252837963122146 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
